# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with
# the latest scraped values (GitHub Actions cryptos-list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.484.20"
$ws.Range("E2").Value = "  +3.57%  "

$ws.Range("D3").Value = "2.759.80"
$ws.Range("E3").Value = "  +4.48%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'116.12"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +2.26%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'333.09"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +2.61%  "

$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("E8").Value = "  +0.00%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.573"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +5.05%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'41.76"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +4.30%  "

$ws.Range("E11").Value = "  +6.55%  "

$ws.Range("E12").Value = "  +2.28%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.129"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +2.20%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'7.65"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +4.45%  "

$ws.Range("D15").Value = "3.192.63"
$ws.Range("E15").Value = "  +4.52%  "

$ws.Range("D16").Value = "2.772.29"
$ws.Range("E16").Value = "  +4.61%  "

$ws.Range("E17").Value = "  +3.18%  "

$ws.Range("D18").Value = "51.533.20"
$ws.Range("E18").Value = "  +3.88%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'3.28"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +9.77%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'13.47"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +4.15%  "

$ws.Range("E21").Value = "  +2.01%  "

$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  +2.92%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'278.15"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +2.88%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'69.67"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.98%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.69"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +5.39%  "

$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("E27").Value = "  +0.02%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'10.16"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -1.93%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'2.23"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +1.38%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'35.02"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -0.73%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'50.08"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +0.92%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = "'5.56"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +1.23%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.0821"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("E35").Value = "  -0.04%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = "'18.96"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -0.61%  "

$ws.Range("E38").Value = "  +1.19%  "

$ws.Range("E39").Value = "  +3.35%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.0352"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +7.63%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'127.15"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +0.14%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").Value = "'23.12"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +2.79%  "

$ws.Range("E43").Value = "  +2.78%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'2.30"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +7.14%  "

$ws.Range("E45").Value = "  +14.58%  "

$ws.Range("D46").Value = "2.088.61"
$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("E47").Value = "  +2.31%  "

$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("E49").Value = "  +5.56%  "

$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("E51").Value = "  +1.39%  "
